# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8085.1333
$ws.Range("I40").Value = 7164.778
$ws.Range("J40").Value = 9465.666999999999
$ws.Range("K40").Value = 7164.778
$ws.Range("L40").Value = 9465.666999999999
$ws.Range("M40").Value = -6989.778
$ws.Range("N40").Value = -9815.666999999999

$ws.Range("H74").Value = 7543
$ws.Range("I74").Value = 7543
$ws.Range("K74").Value = 7543
$ws.Range("M74").Value = -6607

$ws.Range("H77").Value = 7543
$ws.Range("I77").Value = 7543
$ws.Range("K77").Value = 37715
$ws.Range("M77").Value = -33035

$ws.Range("H116").Value = 5259.405
$ws.Range("I116").Value = 6432.636
$ws.Range("J116").Value = 3968.85
$ws.Range("K116").Value = 6432.636
$ws.Range("L116").Value = 3968.85
$ws.Range("M116").Value = -2990.636
$ws.Range("N116").Value = -10852.85

$ws.Range("H136").Value = 137798.17
$ws.Range("J136").Value = 137798.17
$ws.Range("L136").Value = 137798.17
$ws.Range("N136").Value = -147998.17

$ws.Range("H141").Value = 10204.609
$ws.Range("I141").Value = 7035.758
$ws.Range("J141").Value = 23276.125
$ws.Range("K141").Value = 21107.274
$ws.Range("L141").Value = 69828.375
$ws.Range("M141").Value = -15927.274
$ws.Range("N141").Value = -80188.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5797.0864
$ws.Range("I61").Value = 4620.674
$ws.Range("K61").Value = 4620.674
$ws.Range("M61").Value = -4408.674

$ws.Range("H74").Value = 2205.75
$ws.Range("I74").Value = 1924.381
$ws.Range("K74").Value = 1924.381
$ws.Range("M74").Value = -1050.381

$ws.Range("H77").Value = 2205.75
$ws.Range("I77").Value = 1924.381
$ws.Range("K77").Value = 9621.905000000001
$ws.Range("M77").Value = -5253.905000000001

$ws.Range("H136").Value = 5797.0864
$ws.Range("I136").Value = 4620.674
$ws.Range("K136").Value = 13862.022
$ws.Range("M136").Value = -11312.022

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2116.2046
$ws.Range("I20").Value = 1948
$ws.Range("K20").Value = 1948
$ws.Range("M20").Value = -1701

$ws.Range("H22").Value = 349.25
$ws.Range("I22").Value = 349.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 349.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -176.25
$ws.Range("N22").ClearContents()

$ws.Range("H86").Value = 3333
$ws.Range("I86").Value = 3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3333
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2210
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3333
$ws.Range("I89").Value = 3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16665
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11049
$ws.Range("N89").ClearContents()

$ws.Range("H99").Value = 2353.724
$ws.Range("I99").Value = 1774.6666
$ws.Range("K99").Value = 1774.6666
$ws.Range("M99").Value = -276.6666

$ws.Range("H134").Value = 3699.6128
$ws.Range("I134").Value = 2952
$ws.Range("K134").Value = 8856
$ws.Range("M134").Value = -6321

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 9165.076999999999
$ws.Range("J22").Value = 166.25
$ws.Range("L22").Value = 166.25
$ws.Range("N22").Value = -866.25

$ws.Range("H31").Value = 2919.4595
$ws.Range("I31").Value = 2123.7036
$ws.Range("J31").Value = 5068
$ws.Range("K31").Value = 2123.7036
$ws.Range("L31").Value = 5068
$ws.Range("M31").Value = -1828.7036
$ws.Range("N31").Value = -5658

$ws.Range("H34").Value = 2919.4595
$ws.Range("I34").Value = 2123.7036
$ws.Range("J34").Value = 5068
$ws.Range("K34").Value = 2123.7036
$ws.Range("L34").Value = 5068
$ws.Range("M34").Value = -1921.7036
$ws.Range("N34").Value = -5472

$ws.Range("H52").Value = 84999.5
$ws.Range("J52").Value = 84999.5
$ws.Range("L52").Value = 84999.5
$ws.Range("N52").Value = -85587.5

$ws.Range("H58").Value = 2486.1538
$ws.Range("I58").Value = 1949.75
$ws.Range("J58").Value = 3344.4
$ws.Range("K58").Value = 1949.75
$ws.Range("L58").Value = 3344.4
$ws.Range("M58").Value = -1746.75
$ws.Range("N58").Value = -3750.4

$ws.Range("H99").Value = 13989.048
$ws.Range("I99").Value = 9866.799999999999
$ws.Range("J99").Value = 17736.545
$ws.Range("K99").Value = 9866.799999999999
$ws.Range("L99").Value = 17736.545
$ws.Range("M99").Value = -8368.799999999999
$ws.Range("N99").Value = -20732.545

$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()

$ws.Range("H126").Value = 13989.048
$ws.Range("I126").Value = 9866.799999999999
$ws.Range("J126").Value = 17736.545
$ws.Range("K126").Value = 29600.4
$ws.Range("L126").Value = 53209.63499999999
$ws.Range("M126").Value = -27130.4
$ws.Range("N126").Value = -58149.63499999999

$ws.Range("H136").Value = 2486.1538
$ws.Range("I136").Value = 1949.75
$ws.Range("J136").Value = 3344.4
$ws.Range("K136").Value = 5849.25
$ws.Range("L136").Value = 10033.2
$ws.Range("M136").Value = -3299.25
$ws.Range("N136").Value = -15133.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 148.45454
$ws.Range("I61").Value = 132.66667
$ws.Range("K61").Value = 398.00001
$ws.Range("M61").Value = -183.00001

$ws.Range("H68").Value = 5756.3335
$ws.Range("I68").Value = 4998
$ws.Range("J68").Value = 5908
$ws.Range("K68").Value = 14994
$ws.Range("L68").Value = 17724
$ws.Range("M68").Value = -14183
$ws.Range("N68").Value = -19346

$ws.Range("H71").Value = 5756.3335
$ws.Range("I71").Value = 4998
$ws.Range("J71").Value = 5908
$ws.Range("K71").Value = 44982
$ws.Range("L71").Value = 53172
$ws.Range("M71").Value = -40926
$ws.Range("N71").Value = -61284

$ws.Range("H104").Value = 7635.5
$ws.Range("I104").Value = 3328.5
$ws.Range("K104").Value = 9985.5
$ws.Range("M104").Value = -7364.5

$ws.Range("H123").Value = 5357.143

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -2560

$ws.Range("H132").Value = 3424.5833
$ws.Range("I132").Value = 2498
$ws.Range("J132").Value = 3508.818
$ws.Range("K132").Value = 22482
$ws.Range("L132").Value = 31579.362
$ws.Range("M132").Value = -19952
$ws.Range("N132").Value = -36639.362

$ws.Range("H133").Value = 15749.5
$ws.Range("I133").Value = 15333
$ws.Range("J133").Value = 16999
$ws.Range("K133").Value = 45999
$ws.Range("L133").Value = 50997
$ws.Range("M133").Value = -40939
$ws.Range("N133").Value = -61117

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 130289.664
$ws.Range("J42").Value = 130289.664
$ws.Range("L42").Value = 130289.664
$ws.Range("N42").Value = -131259.664

$ws.Range("H102").Value = 8257.375
$ws.Range("I102").Value = 6788.5557
$ws.Range("K102").Value = 6788.5557
$ws.Range("M102").Value = -5166.5557

$ws.Range("H113").Value = 4165
$ws.Range("I113").Value = 3997
$ws.Range("K113").Value = 3997
$ws.Range("M113").Value = -1827

$ws.Range("H115").Value = 130289.664
$ws.Range("J115").Value = 130289.664
$ws.Range("L115").Value = 130289.664
$ws.Range("N115").Value = -132639.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1136
$ws.Range("I22").Value = 1084.7142
$ws.Range("K22").Value = 1084.7142
$ws.Range("M22").Value = -789.7141999999999

$ws.Range("H27").Value = 1136
$ws.Range("I27").Value = 1084.7142
$ws.Range("K27").Value = 1084.7142
$ws.Range("M27").Value = -977.7141999999999

$ws.Range("H40").Value = 5675.4287
$ws.Range("I40").Value = 3007.6667
$ws.Range("K40").Value = 3007.6667
$ws.Range("M40").Value = -2871.6667

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H93").Value = 11119811
$ws.Range("I93").Value = 13340193
$ws.Range("J93").Value = 17901.334
$ws.Range("K93").Value = 13340193
$ws.Range("L93").Value = 17901.334
$ws.Range("M93").Value = -13338945
$ws.Range("N93").Value = -20397.334

$ws.Range("H132").Value = 18752.857
$ws.Range("I132").Value = 20643.084
$ws.Range("J132").Value = 12081.471
$ws.Range("K132").Value = 61929.25199999999
$ws.Range("L132").Value = 36244.413
$ws.Range("M132").Value = -59399.25199999999
$ws.Range("N132").Value = -41304.413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1243.375
$ws.Range("I113").Value = 758.8
$ws.Range("J113").Value = 2051
$ws.Range("K113").Value = 2276.4
$ws.Range("L113").Value = 6153
$ws.Range("M113").Value = -106.3999999999996
$ws.Range("N113").Value = -10493

$ws.Range("H132").Value = 8619.605
$ws.Range("I132").Value = 7007
$ws.Range("K132").Value = 21021
$ws.Range("M132").Value = -18491

$ws.Range("H136").Value = 1969.5735
$ws.Range("I136").Value = 1978.46
$ws.Range("K136").Value = 5935.38
$ws.Range("M136").Value = -3385.38
